$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.281.26'
$ws.Range('E2').Value = '  +6.23%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.789.48'
$ws.Range('E3').Value = '  +22.26%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '621.36'
$ws.Range('E5').Value = '  +8.36%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.23'
$ws.Range('E6').Value = '  +1.81%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.785.75'
$ws.Range('E7').Value = '  +22.21%  '

$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('E9').Value = '  +5.94%  '

$ws.Range('E10').Value = '  +8.86%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.60'
$ws.Range('E11').Value = '  +3.35%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.502'
$ws.Range('E12').Value = '  +7.47%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.57'
$ws.Range('E13').Value = '  +12.38%  '

$ws.Range('E14').Value = '  +6.73%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.420.04'
$ws.Range('E15').Value = '  +22.18%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.794.79'
$ws.Range('E16').Value = '  +22.31%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '71.384.18'
$ws.Range('E17').Value = '  +6.44%  '

$ws.Range('E18').Value = '  +1.65%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.57'
$ws.Range('E19').Value = '  +7.69%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '520.97'
$ws.Range('E20').Value = '  +6.40%  '

$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.97'
$ws.Range('E21').Value = '  +1.35%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.41'
$ws.Range('E22').Value = '  +21.78%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.754'
$ws.Range('E23').Value = '  +9.86%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.01'
$ws.Range('E24').Value = '  +6.66%  '

$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.54'
$ws.Range('E25').Value = '  +11.71%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.61'
$ws.Range('E26').Value = '  +8.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.15'
$ws.Range('E27').Value = '  +9.69%  '

$ws.Range('E28').Value = '  +0.04%  '

$ws.Range('E29').Value = '  +10.97%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.13'
$ws.Range('E30').Value = '  +2.97%  '

$ws.Range('E31').Value = '  +12.28%  '

$ws.Range('E32').Value = '  +20.10%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.30'
$ws.Range('E33').Value = '  +14.90%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.116'
$ws.Range('E34').Value = '  +4.65%  '

$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('E36').Value = '  +12.00%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.16'
$ws.Range('E37').Value = '  +10.56%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.24'
$ws.Range('E38').Value = '  +10.53%  '

$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.343'
$ws.Range('E39').Value = '  +9.92%  '

$ws.Range('E40').Value = '  +9.21%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.74'
$ws.Range('E41').Value = '  +5.31%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '442.88'
$ws.Range('E42').Value = '  +19.95%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.182.47'
$ws.Range('E43').Value = '  +13.63%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.04'
$ws.Range('E44').Value = '  -5.23%  '

$ws.Range('E45').Value = '  +7.85%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.83'
$ws.Range('E46').Value = '  +4.48%  '

$ws.Range('E47').Value = '  +6.25%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.97'
$ws.Range('E48').Value = '  +9.35%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.54'
$ws.Range('E49').Value = '  +3.59%  '

$ws.Range('E50').Value = '  -0.01%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.47'
$ws.Range('E51').Value = '  +7.76%  '
